$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 7) down to the new row 8
# so the new row inherits the same styles (e.g. date format on A, boolean style on G).
$ws.Range("A7:I7").Copy($ws.Range("A8:I8"))

# Now populate the new row's values
$ws.Range("A8").Value = 42650.371990740743
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 10141.870000000001
$ws.Range("D8").Value = 10084.39
$ws.Range("E8").Value = 308.29998799999998
$ws.Range("F8").Value = 306.52999999999997
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = -0.56999999999999995
$ws.Range("I8").Value = $true
